$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Exercise 3 - "Exercise - Create list.py"): fill in Solved / Hint
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "Use join() and split()"

# Row 5 (Exercise 4 - new exercise: "Exercise - Check if number is prime")
$ws.Range("B5").Value = "Exercise - Check if number is prime"
$ws.Range("C5").Value = "Low"
$ws.Range("D5").Value = "Yes"

# Row 6 (Exercise 5 - new exercise: "Exercise: DNA sequencing")
$ws.Range("B6").Value = "Exercise: DNA sequencing"
$ws.Range("C6").Value = "Low"
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "Convert digit except ACTG to ' ' --> join --> Split --> Sort Reverse wrt len"

# Apply highlight fill (theme 9, tint 0.8) to the newly added exercise name cells
$ws.Range("B5:B6").Interior.ThemeColor = 9
$ws.Range("B5:B6").Interior.TintAndShade = 0.8

# Row 6's exercise name uses a plain (non-themed) font color, matching new font style
$ws.Range("B6").Font.Name = "Calibri"
$ws.Range("B6").Font.Size = 11
$ws.Range("B6").Font.ThemeColor = 0
$ws.Range("B6").Font.ColorIndex = 0

# Column E width adjustment (auto-fit style widening to fit new hint text)
$ws.Columns("E").ColumnWidth = 60.1796875

# Update the active selection to reflect where the user left off editing
$ws.Range("E4").Select()
